$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Update refreshed market-data values
$ws.Range("H132").Value = 251831.06
$ws.Range("I132").Value = 312067.9
$ws.Range("K132").Value = 936203.7000000001
$ws.Range("M132").Value = -933673.7000000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Clear market-data columns H:N for rows with no current data
$ws.Range("H117:N120").ClearContents()
$ws.Range("H122:N135").ClearContents()
$ws.Range("H137:N141").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Clear market-data columns H:N for rows with no current data
$ws.Range("H129:N135").ClearContents()
$ws.Range("H137:N141").ClearContents()

# Update refreshed market-data values
$ws.Range("H31").Value = 1971.0209
$ws.Range("I31").Value = 1299.7727
$ws.Range("J31").Value = 2539
$ws.Range("K31").Value = 1299.7727
$ws.Range("L31").Value = 2539
$ws.Range("M31").Value = -1004.7727
$ws.Range("N31").Value = -3129
$ws.Range("H34").Value = 1971.0209
$ws.Range("I34").Value = 1299.7727
$ws.Range("J34").Value = 2539
$ws.Range("K34").Value = 1299.7727
$ws.Range("L34").Value = 2539
$ws.Range("M34").Value = -1097.7727
$ws.Range("N34").Value = -2943

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Clear market-data columns H:N for rows with no current data
$ws.Range("H120:N134").ClearContents()
$ws.Range("H136:N141").ClearContents()

# Update refreshed market-data values
$ws.Range("H68").Value = 1415.34
$ws.Range("I68").Value = 1014.76086
$ws.Range("J68").Value = 1756.5741
$ws.Range("K68").Value = 3044.28258
$ws.Range("L68").Value = 5269.7223
$ws.Range("M68").Value = -2233.28258
$ws.Range("N68").Value = -6891.7223
$ws.Range("H71").Value = 1415.34
$ws.Range("I71").Value = 1014.76086
$ws.Range("J71").Value = 1756.5741
$ws.Range("K71").Value = 9132.847739999999
$ws.Range("L71").Value = 15809.1669
$ws.Range("M71").Value = -5076.847739999999
$ws.Range("N71").Value = -23921.1669
$ws.Range("H80").Value = 1155.7778
$ws.Range("I80").Value = 1049.5
$ws.Range("J80").Value = 1186.1428
$ws.Range("K80").Value = 3148.5
$ws.Range("L80").Value = 3558.4284
$ws.Range("M80").Value = -2212.5
$ws.Range("N80").Value = -5430.428400000001
$ws.Range("H83").Value = 1155.7778
$ws.Range("I83").Value = 1049.5
$ws.Range("J83").Value = 1186.1428
$ws.Range("K83").Value = 9445.5
$ws.Range("L83").Value = 10675.2852
$ws.Range("M83").Value = -4765.5
$ws.Range("N83").Value = -20035.2852

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Update refreshed market-data values
$ws.Range("H132").Value = 4089.5667
$ws.Range("I132").Value = 2514.125
$ws.Range("J132").Value = 5890.0713
$ws.Range("K132").Value = 7542.375
$ws.Range("L132").Value = 17670.2139
$ws.Range("M132").Value = -5012.375
$ws.Range("N132").Value = -22730.2139

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Update refreshed market-data values
$ws.Range("H132").Value = 15627510
$ws.Range("I132").Value = 23811378
$ws.Range("J132").Value = 3763.818
$ws.Range("K132").Value = 71434134
$ws.Range("L132").Value = 11291.454
$ws.Range("M132").Value = -71431604
$ws.Range("N132").Value = -16351.454

